$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: "Fix: Mouse snapped to (0,0) during intro screen" is now complete -> mark Completed? = Yes
# Copy the formatting used by the other "Yes" cells (style with green fill) then set the value.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = "Yes"

# --- Row 24: "Port other UI controls to Zombono UI Engine" also completed, with a
# completion date that still requires polishing (rich text: bold "REQUIRES POLISHING").
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Range("C24").Value = "Yes"

$ws.Range("D24").Value = "05/06/2024 (REQUIRES POLISHING)"
$ws.Range("D24").Characters(13, 18).Font.Bold = $true
$ws.Range("D24").Characters(31, 1).Font.Bold = $false

$ws.Rows(24).RowHeight = 15

# --- Selection moved to A13 (last thing clicked before saving)
$ws.Range("A13").Select() | Out-Null

$excel.CutCopyMode = 0
